# "added futures to model, removed avg_price"
#
# Adds a new "true_values" worksheet (holding the realized/actual values
# used to score the predictions) at the end of the workbook, and updates
# a handful of view-state selections left over from the editing session.
# The "avg_price_all" column header is kept (for layout parity with the
# other prediction sheets) but its data was intentionally left blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new sheet after the last existing sheet, name it, and make
#    it the active tab (this also updates workbookView.activeTab and
#    clears tabSelected from whichever sheet used to be active).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "true_values"

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "dollar_euros"
$ws.Range("C1").Value = "opec_prod"
$ws.Range("D1").Value = "world_demand"
$ws.Range("E1").Value = "min_temp"
$ws.Range("F1").Value = "avg_price_all"
$ws.Range("G1").Value = "s_p_close"
$ws.Range("H1").Value = "usd_fx_index"
$ws.Range("I1").Value = "dow_dji_close"
$ws.Range("J1").Value = "emerging_market_etf"

# ---------------------------------------------------------------------
# 3. Data rows 2-15.
#    Column A = date (copy the date-format from misc_pred so it reuses
#    the existing numFmt "m/d/yyyy" style instead of minting a new one).
# ---------------------------------------------------------------------
$miscPred = $wb.Worksheets.Item("misc_pred")
$miscPred.Range("A2:A15").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122) # xlPasteFormats

$Acol = @(43581, 43582, 43583, 43584, 43585, 43586, 43587, 43588, 43589, 43590, 43591, 43592, 43593, 43594)
$Bcol = @(1.1154999999999999, 1.1156999999999999, 1.115, 1.1184000000000001, 1.1212, 1.1203000000000001, 1.1174999999999999, 1.1207, 1.1207, 1.1200000000000001, 1.1197999999999999, 1.1193, 1.1197999999999999, 1.1220000000000001)
$Ccol = @(30031, 30031, 30031, 30031, 30031, 29876, 29876, 29876, 29876, 29876, 29876, 29876, 29876, 29876)
$Dcol = @(99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2, 99.2)
$Ecol = @(10.3, 10, 9.1, 7.4, 8.3000000000000007, 6.8, 2.9, 0.6, 2.8, 3.4, 2.2999999999999998, 1.7, 0.8, 2.7)
$Gcol = @(2939.88, $null, $null, 2943.03, 2945.83, 2923.73, 2917.52, 2945.64, $null, $null, 2932.47, 2884.05, 2879.42, 2870.72)
$Hcol = @(98.01, $null, $null, 97.86, 97.48, 97.69, 97.83, 97.52, $null, $null, 97.52, 97.63, 97.62, 97.37)
$Icol = @(26543.33, $null, $null, 26554.39, 26592.91, 26430.14, 26307.79, 26504.95, $null, $null, 26438.48, 25965.09, 25967.33, 25828.36)
$Jcol = @(40.03, $null, $null, 39.950000000000003, 39.979999999999997, 39.76, 39.81, 40.33, $null, $null, 39.47, 38.729999999999997, 38.630000000000003, 38.08)

# NB: column F (avg_price_all) is intentionally left empty for every row
# - the commit message is "added futures to model, removed avg_price".

for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $Acol[$i]
    if ($Bcol[$i] -ne $null) { $ws.Cells.Item($r, 2).Value = $Bcol[$i] }
    if ($Ccol[$i] -ne $null) { $ws.Cells.Item($r, 3).Value = $Ccol[$i] }
    if ($Dcol[$i] -ne $null) { $ws.Cells.Item($r, 4).Value = $Dcol[$i] }
    if ($Ecol[$i] -ne $null) { $ws.Cells.Item($r, 5).Value = $Ecol[$i] }
    if ($Gcol[$i] -ne $null) { $ws.Cells.Item($r, 7).Value = $Gcol[$i] }
    if ($Hcol[$i] -ne $null) { $ws.Cells.Item($r, 8).Value = $Hcol[$i] }
    if ($Icol[$i] -ne $null) { $ws.Cells.Item($r, 9).Value = $Icol[$i] }
    if ($Jcol[$i] -ne $null) { $ws.Cells.Item($r, 10).Value = $Jcol[$i] }
}

# G2 (first futures reading) picked up a thousands-format while it was
# being typed in.
$ws.Range("G2").NumberFormat = "#,##0.00"

# O2:O3 carry left-over date formatting from a stray paste (no values).
$miscPred.Range("A2:A3").Copy()
$ws.Range("O2:O3").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Leftover view-state selections on the other sheets from the same
#    editing session.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("OLD").Range("B2:F2").Select()
$wb.Worksheets.Item("weather_pred").Range("A1:A15").Select()
$wb.Worksheets.Item("misc_pred").Range("A1:D15").Select()

# ---------------------------------------------------------------------
# 5. Leave the new sheet active/selected, matching where the author
#    left the cursor.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("R13").Select()
